# Release notes for v0.6.0 — update leaderboard workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename sheets 0.5.2-* -> 0.6.0-*
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)
$wsScores   = $wb.Worksheets.Item(2)
$wsOverview.Name = "0.6.0-Overview"
$wsScores.Name   = "0.6.0-F1-Scores"

# ---------------------------------------------------------------------
# 2. Overview sheet: insert one new data row (for the new "lnn" pipeline)
#    and rewrite every data row in the new, re-sorted order.
# ---------------------------------------------------------------------
$wsOverview.Rows.Item(11).Insert()

# Pipeline, #Wins, #Anomalies, Average F1 Score, Failure Rate
$wsOverview.Cells.Item(2,1).Value  = "aer"
$wsOverview.Cells.Item(2,2).Value  = 11
$wsOverview.Cells.Item(2,3).Value  = 2244
$wsOverview.Cells.Item(2,4).Value  = 0.7295
$wsOverview.Cells.Item(2,5).Value  = 0

$wsOverview.Cells.Item(3,1).Value  = "lstm_dynamic_threshold"
$wsOverview.Cells.Item(3,2).Value  = 8
$wsOverview.Cells.Item(3,3).Value  = 2087
$wsOverview.Cells.Item(3,4).Value  = 0.6278
$wsOverview.Cells.Item(3,5).Value  = 0

$wsOverview.Cells.Item(4,1).Value  = "lstm_autoencoder"
$wsOverview.Cells.Item(4,2).Value  = 7
$wsOverview.Cells.Item(4,3).Value  = 1617
$wsOverview.Cells.Item(4,4).Value  = 0.5486
$wsOverview.Cells.Item(4,5).Value  = 0

$wsOverview.Cells.Item(5,1).Value  = "arima"
$wsOverview.Cells.Item(5,2).ClearContents()
$wsOverview.Cells.Item(5,3).Value  = 2135
$wsOverview.Cells.Item(5,4).Value  = 0.5463
$wsOverview.Cells.Item(5,5).Value  = 0.001

$wsOverview.Cells.Item(6,1).Value  = "vae"
$wsOverview.Cells.Item(6,2).Value  = 6
$wsOverview.Cells.Item(6,3).Value  = 1674
$wsOverview.Cells.Item(6,4).Value  = 0.5375
$wsOverview.Cells.Item(6,5).Value  = 0

$wsOverview.Cells.Item(7,1).Value  = "tadgan"
$wsOverview.Cells.Item(7,2).Value  = 7
$wsOverview.Cells.Item(7,3).Value  = 1816
$wsOverview.Cells.Item(7,4).Value  = 0.536
$wsOverview.Cells.Item(7,5).Value  = 0.0016

$wsOverview.Cells.Item(8,1).Value  = "lnn"
$wsOverview.Cells.Item(8,2).Value  = 7
$wsOverview.Cells.Item(8,3).Value  = 1487
$wsOverview.Cells.Item(8,4).Value  = 0.5351
$wsOverview.Cells.Item(8,5).Value  = 0

$wsOverview.Cells.Item(9,1).Value  = "matrixprofile"
$wsOverview.Cells.Item(9,2).Value  = 5
$wsOverview.Cells.Item(9,3).Value  = 6203
$wsOverview.Cells.Item(9,4).Value  = 0.5268
$wsOverview.Cells.Item(9,5).Value  = 0

$wsOverview.Cells.Item(10,1).Value = "dense_autoencoder"
$wsOverview.Cells.Item(10,2).Value = 7
$wsOverview.Cells.Item(10,3).Value = 964
$wsOverview.Cells.Item(10,4).Value = 0.514
$wsOverview.Cells.Item(10,5).Value = 0

$wsOverview.Cells.Item(11,1).Value = "ganf"
$wsOverview.Cells.Item(11,2).Value = 5
$wsOverview.Cells.Item(11,3).Value = 793
$wsOverview.Cells.Item(11,4).Value = 0.3577
$wsOverview.Cells.Item(11,5).Value = 0.0296

$wsOverview.Cells.Item(12,1).Value = "azure"
$wsOverview.Cells.Item(12,2).Value = 0
$wsOverview.Cells.Item(12,3).Value = 20912
$wsOverview.Cells.Item(12,4).Value = 0.2239
$wsOverview.Cells.Item(12,5).Value = 0

# ---------------------------------------------------------------------
# 3. F1-Scores sheet: insert one new data row (for "lnn") and rewrite
#    every data row's label + values.
# ---------------------------------------------------------------------
$wsScores.Rows.Item(8).Insert()

function Set-ScoreRow($row, $label, $vals) {
    $wsScores.Cells.Item($row, 1).Value = $label
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $wsScores.Cells.Item($row, 2 + $i).Value = $vals[$i]
    }
}

Set-ScoreRow 3  "aer"                    @(0.5867,0.8189,0.4762,0.7988,0.9875,0.892,0.7093,0.7143,0.7407,0.6897,0.7027,0.6377,0.7295,0.1297)
Set-ScoreRow 4  "arima"                  @(0.525,0.4115,0.1533,0.7282,0.8559,0.7972,0.6861,0.3077,0.3824,0.7273,0.4667,0.5143,0.5463,0.2058)
Set-ScoreRow 5  "azure"                  @(0.0512,0.0187,0.0149,0.2796,0.6525,0.7024,0.3442,0.0556,0.1124,0.1626,0.1167,0.1759,0.2239,0.2244)
Set-ScoreRow 6  "dense_autoencoder"      @(0.5588,0.6923,0.2074,0.6667,0.8924,0.0697,0.1008,0.5455,0.7636,0.6,0.5625,0.5085,0.514,0.2476)
Set-ScoreRow 7  "ganf"                   @(0.4615,0.4632,0.1473,0.0857,0.1714,0.0085,0.1525,0.6667,0.5778,0.3077,0.5833,0.6667,0.3577,0.2294)
Set-ScoreRow 8  "lnn"                    @(0.5169,0.6182,0.362,0.6522,0.9381,0.331,0.1908,0.375,0.481,0.7143,0.6667,0.5753,0.5351,0.1943)
Set-ScoreRow 9  "lstm_autoencoder"       @(0.5455,0.6621,0.327,0.5952,0.8667,0.4659,0.2385,0.6667,0.7407,0.5,0.5,0.4746,0.5486,0.1657)
Set-ScoreRow 10 "lstm_dynamic_threshold" @(0.4706,0.726,0.3934,0.7277,0.985,0.744,0.6456,0.4,0.4675,0.7857,0.5854,0.6027,0.6278,0.1697)
Set-ScoreRow 11 "matrixprofile"          @(0.4737,0.4228,0.051,0.5073,0.8975,0.793,0.8251,0.5714,0.44,0.6923,0.3051,0.3429,0.5268,0.2341)
Set-ScoreRow 12 "tadgan"                 @(0.56,0.6049,0.17,0.5775,0.8165,0.4164,0.3404,0.5,0.623,0.8182,0.4516,0.5538,0.536,0.1752)
Set-ScoreRow 13 "vae"                    @(0.4938,0.6131,0.324,0.5924,0.8034,0.4383,0.2303,0.6667,0.6885,0.5833,0.4828,0.5333,0.5375,0.1516)
